$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 154, pushing the existing rows 154-170 down to 155-171
$ws.Rows.Item(154).Insert()

# Populate the newly inserted row 154 with the new weekly record.
# (Columns A,B,C,E,F,G,H,I,N,O,Q,R repeat the same constant values used
# throughout this block; D/J/K/L/M/P carry the new record's data.)
$ws.Cells.Item(154, 1).Value = 9
$ws.Cells.Item(154, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(154, 3).Value = "Metropolitana"
$ws.Cells.Item(154, 4).Value = 44491
$ws.Cells.Item(154, 5).Value = 13
$ws.Cells.Item(154, 6).Value = 300000001
$ws.Cells.Item(154, 7).Value = "Rabanito"
$ws.Cells.Item(154, 8).Value = "Sin especificar"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 8800
$ws.Cells.Item(154, 11).Value = 3000
$ws.Cells.Item(154, 12).Value = 4000
$ws.Cells.Item(154, 13).Value = 3500
$ws.Cells.Item(154, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(154, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(154, 16).Value = 35
$ws.Cells.Item(154, 17).Value = 100
$ws.Cells.Item(154, 18).Value = "Hortaliza"
